$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 630.7778
$ws.Range("I19").Value = 437.5
$ws.Range("J19").Value = 686
$ws.Range("K19").Value = 437.5
$ws.Range("L19").Value = 686
$ws.Range("M19").Value = -262.5
$ws.Range("N19").Value = -1036

$ws.Range("H64").Value = 3103.6667
$ws.Range("I64").Value = 2936.7896
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 2936.7896
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -2688.7896
$ws.Range("N64").Value = -3996

$ws.Range("H67").Value = 3103.6667
$ws.Range("I67").Value = 2936.7896
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 2936.7896
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -2078.7896
$ws.Range("N67").Value = -5216

$ws.Range("H113").Value = 2911.4707
$ws.Range("I113").Value = 3199
$ws.Range("J113").Value = 2791.6667
$ws.Range("K113").Value = 3199
$ws.Range("L113").Value = 2791.6667
$ws.Range("M113").Value = 55
$ws.Range("N113").Value = -9299.6667

$ws.Range("H116").Value = 1753.75
$ws.Range("I116").Value = 1071.6666
$ws.Range("J116").Value = 3800
$ws.Range("K116").Value = 1071.6666
$ws.Range("L116").Value = 3800
$ws.Range("M116").Value = 2370.3334
$ws.Range("N116").Value = -10684

$ws.Range("H132").Value = 7941445.5
$ws.Range("I132").Value = 10004781
$ws.Range("J132").Value = 5541.154
$ws.Range("K132").Value = 30014343
$ws.Range("L132").Value = 16623.462
$ws.Range("M132").Value = -30011813
$ws.Range("N132").Value = -21683.462

$ws.Range("H137").Value = 1211.9584
$ws.Range("I137").Value = 950.5
$ws.Range("J137").Value = 1578
$ws.Range("K137").Value = 2851.5
$ws.Range("L137").Value = 4734
$ws.Range("M137").Value = -301.5
$ws.Range("N137").Value = -9834

$ws.Range("H141").Value = 3169.348
$ws.Range("I141").Value = 2087.7778
$ws.Range("J141").Value = 7063
$ws.Range("K141").Value = 6263.3334
$ws.Range("L141").Value = 21189
$ws.Range("M141").Value = -1083.3334
$ws.Range("N141").Value = -31549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1214.4138
$ws.Range("I45").Value = 1000.9091
$ws.Range("J45").Value = 1885.4286
$ws.Range("K45").Value = 1000.9091
$ws.Range("L45").Value = 1885.4286
$ws.Range("M45").Value = -623.9091
$ws.Range("N45").Value = -2639.4286

$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20636

$ws.Range("H61").Value = 1875.0714
$ws.Range("I61").Value = 1771.75
$ws.Range("J61").Value = 2012.8334
$ws.Range("K61").Value = 1771.75
$ws.Range("L61").Value = 2012.8334
$ws.Range("M61").Value = -1559.75
$ws.Range("N61").Value = -2436.8334

$ws.Range("H132").Value = 1289.66
$ws.Range("I132").Value = 902.125
$ws.Range("K132").Value = 2706.375
$ws.Range("M132").Value = -176.375

$ws.Range("H136").Value = 1875.0714
$ws.Range("I136").Value = 1771.75
$ws.Range("J136").Value = 2012.8334
$ws.Range("K136").Value = 5315.25
$ws.Range("L136").Value = 6038.5002
$ws.Range("M136").Value = -2765.25
$ws.Range("N136").Value = -11138.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 178.86667
$ws.Range("I64").Value = 164
$ws.Range("J64").Value = 191.875
$ws.Range("K64").Value = 164
$ws.Range("L64").Value = 191.875
$ws.Range("M64").Value = 61
$ws.Range("N64").Value = -641.875

$ws.Range("H67").Value = 178.86667
$ws.Range("I67").Value = 164
$ws.Range("J67").Value = 191.875
$ws.Range("K67").Value = 164
$ws.Range("L67").Value = 191.875
$ws.Range("M67").Value = 616
$ws.Range("N67").Value = -1751.875

$ws.Range("H80").Value = 4687.207
$ws.Range("I80").Value = 871.8570999999999
$ws.Range("J80").Value = 8248.200000000001
$ws.Range("K80").Value = 871.8570999999999
$ws.Range("L80").Value = 8248.200000000001
$ws.Range("M80").Value = 126.1429000000001
$ws.Range("N80").Value = -10244.2

$ws.Range("H83").Value = 4687.207
$ws.Range("I83").Value = 871.8570999999999
$ws.Range("J83").Value = 8248.200000000001
$ws.Range("K83").Value = 4359.2855
$ws.Range("L83").Value = 41241
$ws.Range("M83").Value = 632.7145
$ws.Range("N83").Value = -51225

$ws.Range("H94").Value = 534.4828
$ws.Range("I94").Value = 445.82352
$ws.Range("J94").Value = 660.0833
$ws.Range("K94").Value = 445.82352
$ws.Range("L94").Value = 660.0833
$ws.Range("M94").Value = 5.176480000000026
$ws.Range("N94").Value = -1562.0833

$ws.Range("H99").Value = 1707
$ws.Range("I99").Value = 1058.5714
$ws.Range("K99").Value = 1058.5714
$ws.Range("M99").Value = 439.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 37039000
$ws.Range("I99").Value = 111112700
$ws.Range("J99").Value = 2150.6667
$ws.Range("K99").Value = 111112700
$ws.Range("L99").Value = 2150.6667
$ws.Range("M99").Value = -111111202
$ws.Range("N99").Value = -5146.6667

$ws.Range("H126").Value = 37039000
$ws.Range("I126").Value = 111112700
$ws.Range("J126").Value = 2150.6667
$ws.Range("K126").Value = 333338100
$ws.Range("L126").Value = 6452.000100000001
$ws.Range("M126").Value = -333335630
$ws.Range("N126").Value = -11392.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4583852
$ws.Range("I113").Value = 10417128
$ws.Range("J113").Value = 3125533
$ws.Range("K113").Value = 31251384
$ws.Range("L113").Value = 9376599
$ws.Range("M113").Value = -31249214
$ws.Range("N113").Value = -9380939

$ws.Range("H117").Value = 829.9
$ws.Range("I117").Value = 533.3333
$ws.Range("J117").Value = 957
$ws.Range("K117").Value = 1599.9999
$ws.Range("L117").Value = 2871
$ws.Range("M117").Value = 1842.0001
$ws.Range("N117").Value = -9755

$ws.Range("H131").Value = 750.5700000000001
$ws.Range("J131").Value = 777.3587
$ws.Range("L131").Value = 2332.0761
$ws.Range("N131").Value = -12412.0761

$ws.Range("H137").Value = 250002000
$ws.Range("I137").Value = 250002000
$ws.Range("K137").Value = 750006000
$ws.Range("M137").Value = -750000900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1425.091
$ws.Range("I126").Value = 1566
$ws.Range("J126").Value = 1221.5555
$ws.Range("K126").Value = 4698
$ws.Range("L126").Value = 3664.6665
$ws.Range("M126").Value = -2228
$ws.Range("N126").Value = -8604.666499999999

$ws.Range("H132").Value = 5802.2607
$ws.Range("I132").Value = 1172.7
$ws.Range("K132").Value = 3518.1
$ws.Range("M132").Value = -988.1000000000004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 5143
$ws.Range("J64").Value = 5150
$ws.Range("L64").Value = 5150
$ws.Range("N64").Value = -5600

$ws.Range("H67").Value = 5143
$ws.Range("J67").Value = 5150
$ws.Range("L67").Value = 5150
$ws.Range("N67").Value = -6710

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 15041.667
$ws.Range("J63").Value = 15041.667
$ws.Range("L63").Value = 15041.667
$ws.Range("N63").Value = -16289.667

$ws.Range("H66").Value = 15041.667
$ws.Range("J66").Value = 15041.667
$ws.Range("L66").Value = 45125.001
$ws.Range("N66").Value = -51365.001

$ws.Range("H122").Value = 2269.1765
$ws.Range("I122").Value = 1759.5714
$ws.Range("J122").Value = 2625.9
$ws.Range("K122").Value = 5278.7142
$ws.Range("L122").Value = 7877.700000000001
$ws.Range("M122").Value = -2828.7142
$ws.Range("N122").Value = -12777.7

$ws.Range("H126").Value = 785.95654
$ws.Range("I126").Value = 587.3333
$ws.Range("J126").Value = 1501
$ws.Range("K126").Value = 1761.9999
$ws.Range("L126").Value = 4503
$ws.Range("M126").Value = 708.0001
$ws.Range("N126").Value = -9443

$ws.Range("H132").Value = 36649.094
$ws.Range("I132").Value = 47918.316
$ws.Range("J132").Value = 11856.8
$ws.Range("K132").Value = 143754.948
$ws.Range("L132").Value = 35570.39999999999
$ws.Range("M132").Value = -141224.948
$ws.Range("N132").Value = -40630.39999999999
